$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 'Permits Filed for 148–15 Archer Avenue in Jamaica, Queens'
$ws.Range("B2").Value = 'https://newyorkyimby.com/2025/11/permits-filed-for-148-15-archer-avenue-in-jamaica-queens.html'
$ws.Range("C2").Value = 'Permits have been filed for a 22-story affordable housing building at 148–15 Archer Avenue in <a href="https://newyorkyimby.com/neighborhoods/jamaica">Jamaica</a>, Queens. Located at the intersection of 149th Street and Archer Avenue, the lot is near the Sutphin Boulevard–Archer Avenue–JFK Airport subway station, served by the E, J, and Z trains. Larry Grubler of Building 163 SC HDFC is listed as the owner behind the applications.'
$ws.Range("D2").Value = '2025-11-06T11:30:02+00:00'
$ws.Range("E2").Value = 'Thu, 06 Nov 2025 11:30:02 +0000'
$ws.Range("F2").Value = 'YIMBY'
$ws.Range("G2").Value = 'YIMBY - Jamaica'

# content_preview (H2) is blank in the source feed for this entry; leaving
# it unset keeps the cell empty, matching the published row.
